# "Move touch 2 mm right"
#
# S3-pos.xlsx records component placements (Designator, Mid X, Mid Y,
# Rotation, Layer). The touch component move shifted several coordinates
# and, critically, the "J3" placement row was removed from the table
# (its designator no longer exists on the re-positioned part), so every
# row below it (J5, J6, J7, L1, L2, Q1, Q2, Q3, R1, R3..R16, RN1,
# U1..U4, X1, Y1) shifts up by one row, and the last row (56, "Y1")
# disappears from the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BT1 (row 3) moved / rotation reset to 0
$ws.Range("B3").Value = 96.549999999999997
$ws.Range("C3").Value = -70.349999999999994
$ws.Range("D3").Value = 0

# D6 (row 22) moved / rotation reset to 0
$ws.Range("B22").Value = 30.512499999999999
$ws.Range("C22").Value = -60.100000000000001
$ws.Range("D22").Value = 0

# D7 (row 23) moved / rotation reset to 0
$ws.Range("B23").Value = 30.512499999999999
$ws.Range("C23").Value = -63.5
$ws.Range("D23").Value = 0

# The J3 placement (row 26) is gone; deleting it shifts every row below
# (J5..Y1) up by one, so the sheet now ends at row 55 instead of 56.
$ws.Rows.Item(26).Delete()
